$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
        return
    }
    $r.Text = $new
}

# Title
Replace-Text "Astronomy's Eye on Celestial Symphony" "The Essence of Biology Unveiled: Exploring the Realm of Life's Complexities"

# Author
Replace-Text "Isabella Matthews" "Amelia Williams"

# Email - merge "isabella" + "." + "matthews@spatiumobservatory" into one run
Replace-Text "isabella.matthews@spatiumobservatory" "amwilliams97@educonnect"

# Body paragraph 1 (first content paragraph)
Replace-Text "Across the vast panorama of the cosmos, where the reverie of stars enchants the minds of humankind, an intricate symphony of cosmic phenomena awaits our exploration" "Biology, the study of life in all its intricate forms, offers an enthralling journey into the remarkable phenomena that surround us"

Replace-Text " Immerse yourself in this celestial ballet, where galaxies pirouette, planets waltz in orbital harmony, and cosmic winds hum through interstellar voids" " From the swirling patterns of a single cell to the sophisticated symphony of ecosystems, biology unveils the secrets of life's boundless diversity and its inherent unity"

Replace-Text " From the titanic dance of black holes to the radiant birth and death of stars, the universe hums with untold stories. Humans have always looked towards the sky with wonder, curiosity, and a yearning to comprehend the cosmos. In this symphony of celestial wonders, we find a tapestry of knowledge that shapes our understanding of existence, time, and our place in the grand symphony of the universe" " In this realm, we delve into the captivating secrets of living organisms - their structures, intricate processes, and interactions with one another and the wider world"

# Body paragraph 2
Replace-Text "Our journey through this celestial symphony begins with the majestic dance of galaxies, spiraled wonders whirling in their gravitational embrace" "Biology is a captivating mosaic of interconnected concepts, a captivating saga of life's astounding resilience"

Replace-Text " Like graceful ballerinas adoring the cosmic stage, they waltz across vast distances, their intricate choreography revealing the mysteries of gravity and dark matter" " Witnessing the transformation of a caterpillar into a butterfly or unraveling the complex chain of events that allows plants to turn sunlight into energy illuminates the sheer mindboggling capabilities of life"

Replace-Text " As we delve deeper into this celestial waltz, we encounter the captivating ballet of planetary motion, revolving around their stellar partners in an eternal dance of gravitational allure. Intriguing exoplanet discoveries paint a portrait of diverse worlds, ranging from colossal gas giants to terrestrial planets, each harboring secrets waiting to be unraveled" " Biology incites in us an innate quest for answers, unraveling the mysteries of our bodies, understanding the mechanisms behind diseases, and searching for innovative solutions to address pressing ecological issues"

# Body paragraph 3
Replace-Text "Further, this celestial symphony orchestrates a mesmerizing array of celestial spectacles. Meteors streak across the black velvet curtain of the night sky, leaving ephemeral trails of luminescent beauty. Exploding stars erupt in radiant finales, illuminating entire galaxies in their fiery brilliance. Supernovas, akin to cosmic fireworks, forge the elements that shape the universe, while black holes lurk as enigmatic conductors, warbling distorted notes of gravity, and devouring matter, adding a haunting beauty to the cosmic symphony" "As we unravel the intricate tapestry of life, we come face-to-face with questions that have pondered humanity for eons: How did life originate? How do organisms adapt and evolve in response to their surroundings? How can we decipher the genetic language that dictates the symphony of life? Biology provides us with a framework to tackle these perplexing inquiries, inviting us to unlock the secrets of our existence"

# Summary heading paragraph stays "Summary" (unchanged)

# Summary body paragraph
Replace-Text "In the interwoven tapestry of cosmic phenomena, humanity finds a symphony of wonder, knowledge, and boundless exploration" "The study of biology unveils the captivating intricacies of life, from the minuscule world within a single cell to the complex interactions of entire ecosystems"

Replace-Text " This celestial ballet invites introspection, inspiring us to contemplate our place in the universe and our connection to the grand orchestration of existence" " It encapsulates investigations into diverse living organisms, encompassing their structures, functions, and interplay with each other and the environment"

Replace-Text " Astronomy's eye unveils the intricate dance of galaxies, the waltz of planets, and the majestic spectacle of celestial events, painting a portrait of the cosmos that is both awe-inspiring and profound" " Biology inspires us to delve into profound questions about the origins of life, mechanisms of adaptation, and the intricate genetic code"

Replace-Text " Our exploration of this celestial symphony is an odyssey through the vastness of space and time, offering perspectives that transcend earthly boundaries and ignite imaginations" " The journey of understanding biology unveils the essence of our connection with all living things and unravels the fascinating story of our place within the grand tapestry of life on Earth"

# Add a trailing empty paragraph at the end of the document body
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$lastRange.InsertParagraphAfter()

Write-Output "DONE"
